$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.965.78'
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.055.25'
$ws.Range("E3").Value = '  +0.50%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.78'
$ws.Range("E5").Value = '  -1.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.656'
$ws.Range("E6").Value = '  -1.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.79'
$ws.Range("E7").Value = '  -1.90%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '58.64'
$ws.Range("E9").Value = '  -2.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.370'
$ws.Range("E10").Value = '  -3.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0780'
$ws.Range("E11").Value = '  -0.63%  '

$ws.Range("E12").Value = '  +1.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.22'
$ws.Range("E13").Value = '  -3.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.875'
$ws.Range("E14").Value = '  +5.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.352.28'
$ws.Range("E15").Value = '  +0.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.62'
$ws.Range("E16").Value = '  -1.94%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.059.83'
$ws.Range("E17").Value = '  +0.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.908.08'
$ws.Range("E18").Value = '  -0.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.60'
$ws.Range("E19").Value = '  -4.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.16'
$ws.Range("E20").Value = '  -2.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0888'
$ws.Range("E21").Value = '  -0.98%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.40'
$ws.Range("E22").Value = '  +1.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.47'
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("E25").Value = '  +1.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.30'
$ws.Range("E26").Value = '  +10.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.22'
$ws.Range("E27").Value = '  +2.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.33'
$ws.Range("E28").Value = '  -0.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.00'
$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.50'
$ws.Range("E30").Value = '  +15.42%  '

$ws.Range("E31").Value = '  -0.94%  '

$ws.Range("E32").Value = '  -1.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.85'
$ws.Range("E33").Value = '  +6.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0615'
$ws.Range("E34").Value = '  -1.77%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.35'
$ws.Range("E35").Value = '  +6.36%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("E37").Value = '  +4.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0828'
$ws.Range("E38").Value = '  -7.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.96'
$ws.Range("E40").Value = '  -5.16%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0223'
$ws.Range("E41").Value = '  +0.69%  '

$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.96'
$ws.Range("E42").Value = '  -7.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.14'
$ws.Range("E43").Value = '  +1.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.95'
$ws.Range("E44").Value = '  -2.53%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.61'
$ws.Range("E45").Value = '  +0.66%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0943'
$ws.Range("E46").Value = '  -11.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.312.52'
$ws.Range("E47").Value = '  +2.49%  '

$ws.Range("E48").Value = '  -4.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.86'
$ws.Range("E49").Value = '  -1.77%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.75'
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.239.97'
$ws.Range("E51").Value = '  +0.65%  '
